$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above row 319 (all existing rows 319:406 shift down to 320:407)
$ws.Rows.Item(319).Insert()

# Populate the newly inserted row with the new price-report entry
$ws.Cells.Item(319, 1).Value = 4
$ws.Cells.Item(319, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(319, 3).Value = "Los Lagos"
$ws.Cells.Item(319, 4).Value = 44736
$ws.Cells.Item(319, 5).Value = 10
$ws.Cells.Item(319, 6).Value = 100114001
$ws.Cells.Item(319, 7).Value = "Papa"
$ws.Cells.Item(319, 8).Value = "Patagonia"
$ws.Cells.Item(319, 9).Value = "1a (guarda)"
$ws.Cells.Item(319, 10).Value = 600
$ws.Cells.Item(319, 11).Value = 7000
$ws.Cells.Item(319, 12).Value = 8000
$ws.Cells.Item(319, 13).Value = 7500
$ws.Cells.Item(319, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(319, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(319, 16).Value = 300
$ws.Cells.Item(319, 17).Value = 25
$ws.Cells.Item(319, 18).Value = "Hortaliza"
